$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.35335373878479
$ws.Range("B1").Value = 2.37708306312561
$ws.Range("C1").Value = 3.096477508544922
$ws.Range("D1").Value = 3.60206937789917
$ws.Range("E1").Value = 1.897168636322021
